$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing client) is updated: Apellido/DNI/Producto/Total/Cuotas change
$ws.Cells.Item(2, 1).Value = "Nicolas"
$ws.Cells.Item(2, 2).Value = "Davalos"
$ws.Cells.Item(2, 3).Value = "D"
$ws.Cells.Item(2, 4).Value = "quincenal"
$ws.Cells.Item(2, 5).Value = "Tv usado"
$ws.Cells.Item(2, 6).Value = 920000.0
$ws.Cells.Item(2, 7).Value = 0.0
$ws.Cells.Item(2, 8).Value = '[{"numeroCuota":1,"montoOriginal":115000.0,"montoPagado":115000.0,"fechaVencimiento":"2025-07-03","fechaPago":"2025-07-03","isFaltante":false},{"numeroCuota":2,"montoOriginal":115000.0,"montoPagado":0.0,"fechaVencimiento":"2025-07-17","fechaPago":"","isFaltante":false},{"numeroCuota":3,"montoOriginal":115000.0,"montoPagado":0.0,"fechaVencimiento":"2025-07-31","fechaPago":"","isFaltante":false},{"numeroCuota":4,"montoOriginal":115000.0,"montoPagado":0.0,"fechaVencimiento":"2025-08-14","fechaPago":"","isFaltante":false},{"numeroCuota":5,"montoOriginal":115000.0,"montoPagado":0.0,"fechaVencimiento":"2025-08-28","fechaPago":"","isFaltante":false},{"numeroCuota":6,"montoOriginal":115000.0,"montoPagado":0.0,"fechaVencimiento":"2025-09-11","fechaPago":"","isFaltante":false},{"numeroCuota":7,"montoOriginal":115000.0,"montoPagado":0.0,"fechaVencimiento":"2025-09-25","fechaPago":"","isFaltante":false},{"numeroCuota":8,"montoOriginal":115000.0,"montoPagado":0.0,"fechaVencimiento":"2025-10-09","fechaPago":"","isFaltante":false}]'

# Row 3 (new client: Kevincito Chagaray)
$ws.Cells.Item(3, 1).Value = "Kevincito"
$ws.Cells.Item(3, 2).Value = "Chagaray"
$ws.Cells.Item(3, 3).Value = "C"
$ws.Cells.Item(3, 4).Value = "mensual"
$ws.Cells.Item(3, 5).Value = "Celular"
$ws.Cells.Item(3, 6).Value = 270000.0
$ws.Cells.Item(3, 7).Value = 0.0
$ws.Cells.Item(3, 8).Value = '[{"numeroCuota":1,"montoOriginal":90000.0,"montoPagado":90000.0,"fechaVencimiento":"2025-07-03","fechaPago":"2025-07-03","isFaltante":false},{"numeroCuota":2,"montoOriginal":90000.0,"montoPagado":50000.0,"fechaVencimiento":"2025-08-03","fechaPago":"","isFaltante":false},{"numeroCuota":3,"montoOriginal":90000.0,"montoPagado":0.0,"fechaVencimiento":"2025-09-03","fechaPago":"","isFaltante":false}]'

# Row 4 (new client: Kevin Nogueroles)
$ws.Cells.Item(4, 1).Value = "Kevin"
$ws.Cells.Item(4, 2).Value = "Nogueroles"
$ws.Cells.Item(4, 3).Value = "N"
$ws.Cells.Item(4, 4).Value = "quincenal"
$ws.Cells.Item(4, 5).Value = "Tv"
$ws.Cells.Item(4, 6).Value = 1200000.0
$ws.Cells.Item(4, 7).Value = 0.0
$ws.Cells.Item(4, 8).Value = '[{"numeroCuota":1,"montoOriginal":150000.0,"montoPagado":150000.0,"fechaVencimiento":"2025-07-03","fechaPago":"2025-07-03","isFaltante":false},{"numeroCuota":2,"montoOriginal":150000.0,"montoPagado":150000.0,"fechaVencimiento":"2025-07-17","fechaPago":"2025-07-03","isFaltante":false},{"numeroCuota":3,"montoOriginal":150000.0,"montoPagado":100000.0,"fechaVencimiento":"2025-07-31","fechaPago":"","isFaltante":false},{"numeroCuota":4,"montoOriginal":150000.0,"montoPagado":0.0,"fechaVencimiento":"2025-08-14","fechaPago":"","isFaltante":false},{"numeroCuota":5,"montoOriginal":150000.0,"montoPagado":0.0,"fechaVencimiento":"2025-08-28","fechaPago":"","isFaltante":false},{"numeroCuota":6,"montoOriginal":150000.0,"montoPagado":0.0,"fechaVencimiento":"2025-09-11","fechaPago":"","isFaltante":false},{"numeroCuota":7,"montoOriginal":150000.0,"montoPagado":0.0,"fechaVencimiento":"2025-09-25","fechaPago":"","isFaltante":false},{"numeroCuota":8,"montoOriginal":150000.0,"montoPagado":0.0,"fechaVencimiento":"2025-10-09","fechaPago":"","isFaltante":false}]'

# Row 5 (new client: Florencia Nogueroles)
$ws.Cells.Item(5, 1).Value = "Florencia"
$ws.Cells.Item(5, 2).Value = "Nogueroles"
$ws.Cells.Item(5, 3).Value = "N"
$ws.Cells.Item(5, 4).Value = "mensual"
$ws.Cells.Item(5, 5).Value = "Microhondas"
$ws.Cells.Item(5, 6).Value = 200000.0
$ws.Cells.Item(5, 7).Value = 0.0
$ws.Cells.Item(5, 8).Value = '[{"numeroCuota":1,"montoOriginal":50000.0,"montoPagado":0.0,"fechaVencimiento":"2025-07-03","fechaPago":"","isFaltante":false},{"numeroCuota":2,"montoOriginal":50000.0,"montoPagado":0.0,"fechaVencimiento":"2025-08-03","fechaPago":"","isFaltante":false},{"numeroCuota":3,"montoOriginal":50000.0,"montoPagado":0.0,"fechaVencimiento":"2025-09-03","fechaPago":"","isFaltante":false},{"numeroCuota":4,"montoOriginal":50000.0,"montoPagado":0.0,"fechaVencimiento":"2025-10-03","fechaPago":"","isFaltante":false}]'

# Column A and E grew slightly once the new, longer client/product names were
# added ("Kevincito" and "Microhondas"), so re-fit them the way Excel does
# after a paste/entry that changes the longest value in a best-fit column.
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(5).EntireColumn.AutoFit()
